$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.160.76'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '1.640.86'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("E6").Value = '  +2.34%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  +0.03%  '
$cellD10 = $ws.Range("D10")
$styleD10 = $cellD10.Style
$cellD10.Value = "'" + '19.95'
$cellD10.Style = $styleD10
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '1.870.36'
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = '1.640.70'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("E14").Value = '  +0.66%  '
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").Value = '27.164.24'
$ws.Range("E17").Value = '  +1.16%  '
$ws.Range("E18").Value = '  +1.37%  '
$cellD19 = $ws.Range("D19")
$styleD19 = $cellD19.Style
$cellD19.Value = "'" + '217.09'
$cellD19.Style = $styleD19
$ws.Range("E19").Value = '  -1.30%  '
$ws.Range("E20").Value = '  +0.06%  '
$cellD21 = $ws.Range("D21")
$styleD21 = $cellD21.Style
$cellD21.Value = "'" + '6.94'
$cellD21.Style = $styleD21
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("E22").Value = '  +3.33%  '
$cellD23 = $ws.Range("D23")
$styleD23 = $cellD23.Style
$cellD23.Value = "'" + '4.40'
$cellD23.Style = $styleD23
$ws.Range("E23").Value = '  +0.31%  '
$cellD24 = $ws.Range("D24")
$styleD24 = $cellD24.Style
$cellD24.Value = "'" + '9.12'
$cellD24.Style = $styleD24
$ws.Range("E24").Value = '  -0.53%  '
$cellD25 = $ws.Range("D25")
$styleD25 = $cellD25.Style
$cellD25.Value = "'" + '146.59'
$cellD25.Style = $styleD25
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("E27").Value = '  +1.17%  '
$cellD28 = $ws.Range("D28")
$styleD28 = $cellD28.Style
$cellD28.Value = "'" + '0.120'
$cellD28.Style = $styleD28
$ws.Range("E28").Value = '  +0.16%  '
$cellD29 = $ws.Range("D29")
$styleD29 = $cellD29.Style
$cellD29.Value = "'" + '15.68'
$cellD29.Style = $styleD29
$ws.Range("E29").Value = '  -0.56%  '
$cellD30 = $ws.Range("D30")
$styleD30 = $cellD30.Style
$cellD30.Value = "'" + '0.0509'
$cellD30.Style = $styleD30
$ws.Range("E30").Value = '  +1.14%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +1.47%  '
$cellD33 = $ws.Range("D33")
$styleD33 = $cellD33.Style
$cellD33.Value = "'" + '3.01'
$cellD33.Style = $styleD33
$ws.Range("E33").Value = '  +0.47%  '
$ws.Range("D34").Value = '1.304.02'
$ws.Range("E34").Value = '  +3.89%  '
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("E36").Value = '  +1.54%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("E38").Value = '  +2.79%  '
$cellD39 = $ws.Range("D39")
$styleD39 = $cellD39.Style
$cellD39.Value = "'" + '0.859'
$cellD39.Style = $styleD39
$ws.Range("E39").Value = '  +3.09%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("E42").Value = '  +5.62%  '
$cellD43 = $ws.Range("D43")
$styleD43 = $cellD43.Style
$cellD43.Value = "'" + '5.29'
$cellD43.Style = $styleD43
$ws.Range("E43").Value = '  -1.99%  '
$ws.Range("D44").Value = '1.780.39'
$ws.Range("E44").Value = '  +0.20%  '
$cellD45 = $ws.Range("D45")
$styleD45 = $cellD45.Style
$cellD45.Value = "'" + '61.78'
$cellD45.Style = $styleD45
$ws.Range("E45").Value = '  +0.20%  '
$cellD46 = $ws.Range("D46")
$styleD46 = $cellD46.Style
$cellD46.Value = "'" + '91.89'
$cellD46.Style = $styleD46
$ws.Range("E46").Value = '  +0.41%  '
$ws.Range("E47").Value = '  +1.83%  '
$ws.Range("D48").Value = '0.0₆0107'
$ws.Range("E48").Value = '  +1.47%  '
$ws.Range("E49").Value = '  -0.29%  '
$cellD50 = $ws.Range("D50")
$styleD50 = $cellD50.Style
$cellD50.Value = "'" + '7.65'
$cellD50.Style = $styleD50
$ws.Range("E50").Value = '  +0.68%  '
$cellD51 = $ws.Range("D51")
$styleD51 = $cellD51.Style
$cellD51.Value = "'" + '0.0962'
$cellD51.Style = $styleD51
$ws.Range("E51").Value = '  +0.04%  '
